# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-calculation columns (H:N)
# across all 8 job sheets, per the scheduled-runner market-data refresh.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 354.85715
$ws.Range("I2").Value = 339
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 339
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = -226
$ws.Range("N2").Value = -676
$ws.Range("H33").Value = 41856.11
$ws.Range("I33").Value = 50200.734
$ws.Range("J33").Value = 133
$ws.Range("K33").Value = 50200.734
$ws.Range("L33").Value = 133
$ws.Range("M33").Value = -49971.734
$ws.Range("N33").Value = -591
$ws.Range("H40").Value = 2472.6365
$ws.Range("J40").Value = 2562.375
$ws.Range("L40").Value = 2562.375
$ws.Range("N40").Value = -2912.375
$ws.Range("H51").Value = 7444.1113
$ws.Range("I51").Value = 4999.6665
$ws.Range("J51").Value = 8666.333000000001
$ws.Range("K51").Value = 4999.6665
$ws.Range("L51").Value = 8666.333000000001
$ws.Range("M51").Value = -4515.6665
$ws.Range("N51").Value = -9634.333000000001
$ws.Range("H70").Value = 1763.75
$ws.Range("I70").Value = 1727.5
$ws.Range("J70").Value = 1800
$ws.Range("K70").Value = 5182.5
$ws.Range("L70").Value = 5400
$ws.Range("M70").Value = -4912.5
$ws.Range("N70").Value = -5940
$ws.Range("H73").Value = 1763.75
$ws.Range("I73").Value = 1727.5
$ws.Range("J73").Value = 1800
$ws.Range("K73").Value = 5182.5
$ws.Range("L73").Value = 5400
$ws.Range("M73").Value = -4246.5
$ws.Range("N73").Value = -7272
$ws.Range("H80").Value = 12305.454
$ws.Range("J80").Value = 3793.6
$ws.Range("L80").Value = 11380.8
$ws.Range("N80").Value = -13376.8
$ws.Range("H83").Value = 12305.454
$ws.Range("J83").Value = 3793.6
$ws.Range("L83").Value = 34142.4
$ws.Range("N83").Value = -44126.4
$ws.Range("H86").Value = 19548.75
$ws.Range("I86").Value = 9478.1
$ws.Range("J86").Value = 36333.168
$ws.Range("K86").Value = 9478.1
$ws.Range("L86").Value = 36333.168
$ws.Range("M86").Value = -8355.1
$ws.Range("N86").Value = -38579.168
$ws.Range("H89").Value = 19548.75
$ws.Range("I89").Value = 9478.1
$ws.Range("J89").Value = 36333.168
$ws.Range("K89").Value = 47390.5
$ws.Range("L89").Value = 181665.84
$ws.Range("M89").Value = -41774.5
$ws.Range("N89").Value = -192897.84
$ws.Range("H125").Value = 3000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H127").Value = 1017.8461
$ws.Range("I127").Value = 1084.5834
$ws.Range("J127").Value = 217
$ws.Range("K127").Value = 3253.7502
$ws.Range("L127").Value = 651
$ws.Range("M127").Value = 1706.2498
$ws.Range("N127").Value = -10571
$ws.Range("H129").Value = 1288.4546
$ws.Range("I129").Value = 884.75
$ws.Range("J129").Value = 1519.1428
$ws.Range("K129").Value = 2654.25
$ws.Range("L129").Value = 4557.428400000001
$ws.Range("M129").Value = 2345.75
$ws.Range("N129").Value = -14557.4284
$ws.Range("H138").Value = 2800.4717
$ws.Range("I138").Value = 1307.5883
$ws.Range("J138").Value = 3505.4443
$ws.Range("K138").Value = 3922.7649
$ws.Range("L138").Value = 10516.3329
$ws.Range("M138").Value = 1217.2351
$ws.Range("N138").Value = -20796.3329
$ws.Range("H141").Value = 3167.7058
$ws.Range("I141").Value = 2889.4614
$ws.Range("J141").Value = 4072
$ws.Range("K141").Value = 8668.3842
$ws.Range("L141").Value = 12216
$ws.Range("M141").Value = -3488.3842
$ws.Range("N141").Value = -22576

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 25001906
$ws.Range("I132").Value = 31251466
$ws.Range("J132").Value = 3668.75
$ws.Range("K132").Value = 93754398
$ws.Range("L132").Value = 11006.25
$ws.Range("M132").Value = -93751868
$ws.Range("N132").Value = -16066.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 6500
$ws.Range("J34").Value = 6500
$ws.Range("L34").Value = 6500
$ws.Range("N34").Value = -6728
$ws.Range("H94").Value = 3894.1904
$ws.Range("I94").Value = 3634
$ws.Range("K94").Value = 3634
$ws.Range("M94").Value = -3183
$ws.Range("H99").Value = 998.1667
$ws.Range("I99").Value = 918.6667
$ws.Range("J99").Value = 1236.6666
$ws.Range("K99").Value = 918.6667
$ws.Range("L99").Value = 1236.6666
$ws.Range("M99").Value = 579.3333
$ws.Range("N99").Value = -4232.6666

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2982.5293
$ws.Range("I31").Value = 1949.75
$ws.Range("K31").Value = 1949.75
$ws.Range("M31").Value = -1654.75
$ws.Range("H34").Value = 2982.5293
$ws.Range("I34").Value = 1949.75
$ws.Range("K34").Value = 1949.75
$ws.Range("M34").Value = -1747.75
$ws.Range("H58").Value = 1600.091
$ws.Range("I58").Value = 1255.7778
$ws.Range("K58").Value = 1255.7778
$ws.Range("M58").Value = -1052.7778
$ws.Range("H62").Value = 200007200
$ws.Range("I62").Value = 9000
$ws.Range("J62").Value = 500004500
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 500004500
$ws.Range("M62").Value = -8376
$ws.Range("N62").Value = -500005748
$ws.Range("H65").Value = 200007200
$ws.Range("I65").Value = 9000
$ws.Range("J65").Value = 500004500
$ws.Range("K65").Value = 45000
$ws.Range("L65").Value = 2500022500
$ws.Range("M65").Value = -41880
$ws.Range("N65").Value = -2500028740
$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71498
$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -217488
$ws.Range("H105").Value = 1744.5
$ws.Range("I105").Value = 1744.5
$ws.Range("K105").Value = 1744.5
$ws.Range("M105").Value = 2.5
$ws.Range("H136").Value = 1600.091
$ws.Range("I136").Value = 1255.7778
$ws.Range("K136").Value = 3767.3334
$ws.Range("M136").Value = -1217.3334

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2566.2964
$ws.Range("I136").Value = 2881.6667
$ws.Range("J136").Value = 2476.1904
$ws.Range("K136").Value = 8645.000100000001
$ws.Range("L136").Value = 7428.5712
$ws.Range("M136").Value = -3545.000100000001
$ws.Range("N136").Value = -17628.5712
$ws.Range("H137").Value = 1442.4445
$ws.Range("I137").Value = 926
$ws.Range("J137").Value = 3250
$ws.Range("K137").Value = 2778
$ws.Range("L137").Value = 9750
$ws.Range("M137").Value = 2322
$ws.Range("N137").Value = -19950

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 58029.5
$ws.Range("J15").Value = 58029.5
$ws.Range("L15").Value = 58029.5
$ws.Range("N15").Value = -58605.5
$ws.Range("H51").Value = 80000
$ws.Range("J51").Value = 80000
$ws.Range("L51").Value = 80000
$ws.Range("N51").Value = -81018
$ws.Range("H80").Value = 3232.5386
$ws.Range("I80").Value = 2922
$ws.Range("J80").Value = 3426.625
$ws.Range("K80").Value = 2922
$ws.Range("L80").Value = 3426.625
$ws.Range("M80").Value = -1924
$ws.Range("N80").Value = -5422.625
$ws.Range("H81").Value = 58029.5
$ws.Range("J81").Value = 58029.5
$ws.Range("L81").Value = 58029.5
$ws.Range("N81").Value = -60025.5
$ws.Range("H83").Value = 3232.5386
$ws.Range("I83").Value = 2922
$ws.Range("J83").Value = 3426.625
$ws.Range("K83").Value = 14610
$ws.Range("L83").Value = 17133.125
$ws.Range("M83").Value = -9618
$ws.Range("N83").Value = -27117.125
$ws.Range("H84").Value = 58029.5
$ws.Range("J84").Value = 58029.5
$ws.Range("L84").Value = 174088.5
$ws.Range("N84").Value = -184072.5
$ws.Range("H113").Value = 1078
$ws.Range("I113").Value = 1078
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1078
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1092
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2568.8823
$ws.Range("I122").Value = 2206.2307
$ws.Range("J122").Value = 3747.5
$ws.Range("K122").Value = 6618.6921
$ws.Range("L122").Value = 11242.5
$ws.Range("M122").Value = -4168.6921
$ws.Range("N122").Value = -16142.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1708.7084
$ws.Range("I46").Value = 934.63635
$ws.Range("K46").Value = 934.63635
$ws.Range("M46").Value = -746.63635
$ws.Range("H93").Value = 1281
$ws.Range("I93").Value = 1209.1
$ws.Range("K93").Value = 1209.1
$ws.Range("M93").Value = 38.90000000000009
$ws.Range("H132").Value = 3693.4119
$ws.Range("I132").Value = 3461.111
$ws.Range("K132").Value = 10383.333
$ws.Range("M132").Value = -7853.332999999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 7520.75
$ws.Range("I58").Value = 7520.75
$ws.Range("K58").Value = 7520.75
$ws.Range("M58").Value = -7212.75
$ws.Range("H106").Value = 22666.334
$ws.Range("I106").Value = 21499.5
$ws.Range("K106").Value = 21499.5
$ws.Range("M106").Value = -20237.5
$ws.Range("H132").Value = 3305.7556
$ws.Range("I132").Value = 2824.543
$ws.Range("K132").Value = 8473.629000000001
$ws.Range("M132").Value = -5943.629000000001
$ws.Range("H136").Value = 5645.273
$ws.Range("I136").Value = 2299.8572
$ws.Range("K136").Value = 6899.571599999999
$ws.Range("M136").Value = -4349.571599999999
